$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: remove the centered style (revert to default/"Normal" style), value unchanged
$ws.Range("B1").Style = "Normal"
$ws.Range("B1").Value = 50

# B2:E4: remove centered style and write the new gain-table values
$ws.Range("B2:E4").Style = "Normal"

$ws.Range("B2").Value = 3.5362774
$ws.Range("C2").Value = 10.028527
$ws.Range("D2").Value = 0.003
$ws.Range("E2").Value = 0.003

$ws.Range("B3").Value = 0.286032
$ws.Range("C3").Value = 0.118836
$ws.Range("D3").Value = 0.15
$ws.Range("E3").Value = 0.15

$ws.Range("B4").Value = 0.0572064
$ws.Range("C4").Value = 0.023767
$ws.Range("D4").Value = 0.069
$ws.Range("E4").Value = 0.069

# Set column B to the standard default column width (renders as width 10 in the sheet XML)
$ws.Range("B1").ColumnWidth = 9.140625

# Move the active selection to D8
$ws.Range("D8").Select()
